$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("2025").Range("A2").Value = 57
$wb.Worksheets.Item("2030").Range("A2").Value = 195
$wb.Worksheets.Item("2040").Range("A2").Value = 355
$wb.Worksheets.Item("2045").Range("A2").Value = 355
$wb.Worksheets.Item("2050").Range("A2").Value = 355
